$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.748.41'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '1.732.19'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9973'
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.62'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9978'
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4928'
$ws.Range("E7").Value = '  +1.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2627'
$ws.Range("E8").Value = '  +0.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06226'
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D10").Value = '1.725.97'
$ws.Range("E10").Value = '  -0.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '15.93'
$ws.Range("E11").Value = '  +2.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06999'
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6133'
$ws.Range("E13").Value = '  +2.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.505'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '26.532.62'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9972'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007250'
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.45'
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").Value = '1.948.70'
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.496'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.577'
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.117'
$ws.Range("E24").Value = '  -1.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.28'
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.34'
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.779'
$ws.Range("E27").Value = '  +3.26%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.387'
$ws.Range("E28").Value = '  -1.63%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.53'
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.938'
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07994'
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.677'
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04485'
$ws.Range("E33").Value = '  -0.84%  '
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.004'
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6246'
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9360'
$ws.Range("E37").Value = '  +3.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.055'
$ws.Range("E38").Value = '  +3.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.421'
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9977'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01516'
$ws.Range("E41").Value = '  +1.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.620'
$ws.Range("E42").Value = '  +3.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.48'
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3867'
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.931'
$ws.Range("E45").Value = '  +3.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1160'
$ws.Range("E46").Value = '  +0.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05383'
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.919'
$ws.Range("E48").Value = '  +2.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.29'
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.83'
$ws.Range("E50").Value = '  +1.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.237'
$ws.Range("E51").Value = '  -1.51%  '
